# Append a new effort-log entry (row 45) to the "effort" worksheet,
# mirroring the existing rows: date in column A, effort hours in column B,
# description text in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").Value = 41233
$ws.Range("B45").Value = 2.5
$ws.Range("D45").Value = "New test case tc09 put to operation and completed. Fix of makefile. Documentation. Export for Andreas"

# Match the workbook's recorded selection after the edit.
$ws.Range("D45").Select()
